$wb = $excel.ActiveWorkbook

# ---- Sheet "Summary" ----
$ws1 = $wb.Worksheets.Item("Summary")

$ws1.Cells.Item(2,1).Value = 'claude-3-haiku-20240307'
$ws1.Cells.Item(2,2).Value = 444.4
$ws1.Cells.Item(2,3).Value = 234
$ws1.Cells.Item(2,4).Value = 0.000404
$ws1.Cells.Item(2,5).Value = 2.38
$ws1.Cells.Item(2,6).Value = 30.11
$ws1.Cells.Item(2,7).Value = 5

$ws1.Cells.Item(3,1).Value = 'claude-sonnet-4-20250514'
$ws1.Cells.Item(3,2).Value = 444.4
$ws1.Cells.Item(3,3).Value = 244
$ws1.Cells.Item(3,4).Value = 0.004993
$ws1.Cells.Item(3,5).Value = 29.46
$ws1.Cells.Item(3,6).Value = 372.49
$ws1.Cells.Item(3,7).Value = 5

$ws1.Cells.Item(4,1).Value = 'gpt-4o-mini'
$ws1.Cells.Item(4,2).Value = 406.6
$ws1.Cells.Item(4,3).Value = 192.6
$ws1.Cells.Item(4,4).Value = 0.000177
$ws1.Cells.Item(4,5).Value = 1.04
$ws1.Cells.Item(4,6).Value = 13.17
$ws1.Cells.Item(4,7).Value = 5

$ws1.Cells.Item(5,1).Value = 'gpt-4o'
$ws1.Cells.Item(5,2).Value = 406.6
$ws1.Cells.Item(5,3).Value = 189.6
$ws1.Cells.Item(5,4).Value = 0.002912
$ws1.Cells.Item(5,5).Value = 17.18
$ws1.Cells.Item(5,6).Value = 217.27
$ws1.Cells.Item(5,7).Value = 5

# ---- Sheet "Details" ----
$ws2 = $wb.Worksheets.Item("Details")

$ws2.Cells.Item(2,1).Value = 47040
$ws2.Cells.Item(2,2).Value = 'Mind, delusions, eternity, that he was in'
$ws2.Cells.Item(2,3).Value = 'claude-3-haiku-20240307'
$ws2.Cells.Item(2,4).Value = 'a delusional belief that one is in eternity'
$ws2.Cells.Item(2,5).Value = 445
$ws2.Cells.Item(2,6).Value = 231
$ws2.Cells.Item(2,7).Value = 0.0004

$ws2.Cells.Item(3,1).Value = 47040
$ws2.Cells.Item(3,2).Value = 'Mind, delusions, eternity, that he was in'
$ws2.Cells.Item(3,3).Value = 'claude-sonnet-4-20250514'
$ws2.Cells.Item(3,4).Value = 'delusion of being in eternity'
$ws2.Cells.Item(3,5).Value = 445
$ws2.Cells.Item(3,6).Value = 235
$ws2.Cells.Item(3,7).Value = 0.00486

$ws2.Cells.Item(4,1).Value = 47040
$ws2.Cells.Item(4,2).Value = 'Mind, delusions, eternity, that he was in'
$ws2.Cells.Item(4,3).Value = 'gpt-4o-mini'
$ws2.Cells.Item(4,4).Value = 'delusion of being in eternity'
$ws2.Cells.Item(4,5).Value = 407
$ws2.Cells.Item(4,6).Value = 187
$ws2.Cells.Item(4,7).Value = 0.00017325

$ws2.Cells.Item(5,1).Value = 47040
$ws2.Cells.Item(5,2).Value = 'Mind, delusions, eternity, that he was in'
$ws2.Cells.Item(5,3).Value = 'gpt-4o'
$ws2.Cells.Item(5,4).Value = 'belief of being in eternity'
$ws2.Cells.Item(5,5).Value = 407
$ws2.Cells.Item(5,6).Value = 167
$ws2.Cells.Item(5,7).Value = 0.0026875

$ws2.Cells.Item(6,1).Value = 46193
$ws2.Cells.Item(6,2).Value = 'Mind, anxiety, forenoon'
$ws2.Cells.Item(6,3).Value = 'claude-3-haiku-20240307'
$ws2.Cells.Item(6,4).Value = 'Tendency to experience anxiety during the forenoon'
$ws2.Cells.Item(6,5).Value = 440
$ws2.Cells.Item(6,6).Value = 193
$ws2.Cells.Item(6,7).Value = 0.00035125

$ws2.Cells.Item(7,1).Value = 46193
$ws2.Cells.Item(7,2).Value = 'Mind, anxiety, forenoon'
$ws2.Cells.Item(7,3).Value = 'claude-sonnet-4-20250514'
$ws2.Cells.Item(7,4).Value = 'anxiety occurring in the forenoon (late morning hours)'
$ws2.Cells.Item(7,5).Value = 440
$ws2.Cells.Item(7,6).Value = 221
$ws2.Cells.Item(7,7).Value = 0.004635

$ws2.Cells.Item(8,1).Value = 46193
$ws2.Cells.Item(8,2).Value = 'Mind, anxiety, forenoon'
$ws2.Cells.Item(8,3).Value = 'gpt-4o-mini'
$ws2.Cells.Item(8,4).Value = 'anxiety in the forenoon'
$ws2.Cells.Item(8,5).Value = 403
$ws2.Cells.Item(8,6).Value = 179
$ws2.Cells.Item(8,7).Value = 0.00016785

$ws2.Cells.Item(9,1).Value = 46193
$ws2.Cells.Item(9,2).Value = 'Mind, anxiety, forenoon'
$ws2.Cells.Item(9,3).Value = 'gpt-4o'
$ws2.Cells.Item(9,4).Value = 'anxiety in the late morning'
$ws2.Cells.Item(9,5).Value = 403
$ws2.Cells.Item(9,6).Value = 183
$ws2.Cells.Item(9,7).Value = 0.0028375

$ws2.Cells.Item(10,1).Value = 49200
$ws2.Cells.Item(10,2).Value = 'Mind, restlessness, menses, during'
$ws2.Cells.Item(10,3).Value = 'claude-3-haiku-20240307'
$ws2.Cells.Item(10,4).Value = 'Restlessness during menstruation'
$ws2.Cells.Item(10,5).Value = 443
$ws2.Cells.Item(10,6).Value = 258
$ws2.Cells.Item(10,7).Value = 0.00043325

$ws2.Cells.Item(11,1).Value = 49200
$ws2.Cells.Item(11,2).Value = 'Mind, restlessness, menses, during'
$ws2.Cells.Item(11,3).Value = 'claude-sonnet-4-20250514'
$ws2.Cells.Item(11,4).Value = 'restlessness during menstruation'
$ws2.Cells.Item(11,5).Value = 443
$ws2.Cells.Item(11,6).Value = 288
$ws2.Cells.Item(11,7).Value = 0.005649

$ws2.Cells.Item(12,1).Value = 49200
$ws2.Cells.Item(12,2).Value = 'Mind, restlessness, menses, during'
$ws2.Cells.Item(12,3).Value = 'gpt-4o-mini'
$ws2.Cells.Item(12,4).Value = 'restlessness during menstruation'
$ws2.Cells.Item(12,5).Value = 405
$ws2.Cells.Item(12,6).Value = 196
$ws2.Cells.Item(12,7).Value = 0.00017835

$ws2.Cells.Item(13,1).Value = 49200
$ws2.Cells.Item(13,2).Value = 'Mind, restlessness, menses, during'
$ws2.Cells.Item(13,3).Value = 'gpt-4o'
$ws2.Cells.Item(13,4).Value = 'restlessness during menstruation'
$ws2.Cells.Item(13,5).Value = 405
$ws2.Cells.Item(13,6).Value = 197
$ws2.Cells.Item(13,7).Value = 0.0029825

$ws2.Cells.Item(14,1).Value = 48718
$ws2.Cells.Item(14,2).Value = 'Mind, kleptomania, steals money'
$ws2.Cells.Item(14,3).Value = 'claude-3-haiku-20240307'
$ws2.Cells.Item(14,4).Value = 'a tendency to steal money'
$ws2.Cells.Item(14,5).Value = 443
$ws2.Cells.Item(14,6).Value = 205
$ws2.Cells.Item(14,7).Value = 0.000367

$ws2.Cells.Item(15,1).Value = 48718
$ws2.Cells.Item(15,2).Value = 'Mind, kleptomania, steals money'
$ws2.Cells.Item(15,3).Value = 'claude-sonnet-4-20250514'
$ws2.Cells.Item(15,4).Value = 'compulsive stealing of money'
$ws2.Cells.Item(15,5).Value = 443
$ws2.Cells.Item(15,6).Value = 212
$ws2.Cells.Item(15,7).Value = 0.004509

$ws2.Cells.Item(16,1).Value = 48718
$ws2.Cells.Item(16,2).Value = 'Mind, kleptomania, steals money'
$ws2.Cells.Item(16,3).Value = 'gpt-4o-mini'
$ws2.Cells.Item(16,4).Value = 'compulsive stealing of money'
$ws2.Cells.Item(16,5).Value = 405
$ws2.Cells.Item(16,6).Value = 178
$ws2.Cells.Item(16,7).Value = 0.00016755

$ws2.Cells.Item(17,1).Value = 48718
$ws2.Cells.Item(17,2).Value = 'Mind, kleptomania, steals money'
$ws2.Cells.Item(17,3).Value = 'gpt-4o'
$ws2.Cells.Item(17,4).Value = 'a compulsion to steal money'
$ws2.Cells.Item(17,5).Value = 405
$ws2.Cells.Item(17,6).Value = 171
$ws2.Cells.Item(17,7).Value = 0.0027225

$ws2.Cells.Item(18,1).Value = 49192
$ws2.Cells.Item(18,2).Value = 'Mind, restlessness, internal, as if would beat about herself with hands and feet'
$ws2.Cells.Item(18,3).Value = 'claude-3-haiku-20240307'
$ws2.Cells.Item(18,4).Value = 'Restless inner agitation, as if needing to physically move about'
$ws2.Cells.Item(18,5).Value = 451
$ws2.Cells.Item(18,6).Value = 283
$ws2.Cells.Item(18,7).Value = 0.0004665

$ws2.Cells.Item(19,1).Value = 49192
$ws2.Cells.Item(19,2).Value = 'Mind, restlessness, internal, as if would beat about herself with hands and feet'
$ws2.Cells.Item(19,3).Value = 'claude-sonnet-4-20250514'
$ws2.Cells.Item(19,4).Value = 'internal restlessness with urge to thrash about with hands and feet'
$ws2.Cells.Item(19,5).Value = 451
$ws2.Cells.Item(19,6).Value = 264
$ws2.Cells.Item(19,7).Value = 0.005313

$ws2.Cells.Item(20,1).Value = 49192
$ws2.Cells.Item(20,2).Value = 'Mind, restlessness, internal, as if would beat about herself with hands and feet'
$ws2.Cells.Item(20,3).Value = 'gpt-4o-mini'
$ws2.Cells.Item(20,4).Value = 'inner restlessness, feeling a need to thrash about with hands and feet'
$ws2.Cells.Item(20,5).Value = 413
$ws2.Cells.Item(20,6).Value = 223
$ws2.Cells.Item(20,7).Value = 0.00019575

$ws2.Cells.Item(21,1).Value = 49192
$ws2.Cells.Item(21,2).Value = 'Mind, restlessness, internal, as if would beat about herself with hands and feet'
$ws2.Cells.Item(21,3).Value = 'gpt-4o'
$ws2.Cells.Item(21,4).Value = 'feeling an internal restlessness that makes them want to thrash about'
$ws2.Cells.Item(21,5).Value = 413
$ws2.Cells.Item(21,6).Value = 230
$ws2.Cells.Item(21,7).Value = 0.0033325

"edit complete"
